$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.020349605984689
$ws.Range("D2").Value = 1.022174629359064
$ws.Range("E2").Value = 1.029836972115252
$ws.Range("F2").Value = 1.037336229550243
$ws.Range("J2").Value = 1.02554721649973
$ws.Range("K2").Value = 1.025009786906153
$ws.Range("L2").Value = 1.032649726007894
$ws.Range("M2").Value = 1.040127404791578
$ws.Range("N2").Value = 1.012499121358909

# Row 3
$ws.Range("C3").Value = 1.021463424500027
$ws.Range("D3").Value = 1.023144942643183
$ws.Range("E3").Value = 1.030864681332056
$ws.Range("F3").Value = 1.038536172828092
$ws.Range("J3").Value = 1.026297326841434
$ws.Range("K3").Value = 1.025786262692901
$ws.Range("L3").Value = 1.03348507388885
$ws.Range("M3").Value = 1.041136101883477
$ws.Range("N3").Value = 1.012756310177493

# Row 4
$ws.Range("C4").Value = 1.022184695731987
$ws.Range("D4").Value = 1.023773588947807
$ws.Range("E4").Value = 1.031530533551314
$ws.Range("F4").Value = 1.039313808765389
$ws.Range("J4").Value = 1.026782745325733
$ws.Range("K4").Value = 1.026288871349019
$ws.Range("L4").Value = 1.034025853867295
$ws.Range("M4").Value = 1.04178942196542
$ws.Range("N4").Value = 1.012922514777585

# Row 5
$ws.Range("C5").Value = 1.022488051932126
$ws.Range("D5").Value = 1.024038060559988
$ws.Range("E5").Value = 1.031810662680135
$ws.Range("F5").Value = 1.03964101271398
$ws.Range("J5").Value = 1.026986826549574
$ws.Range("K5").Value = 1.026500210362561
$ws.Range("L5").Value = 1.034253258650684
$ws.Range("M5").Value = 1.042064227717361
$ws.Range("N5").Value = 1.012992335863613

# Row 6
$ws.Range("C6").Value = 1.022538994635924
$ws.Range("D6").Value = 1.024082477569908
$ws.Range("E6").Value = 1.031857709606912
$ws.Range("F6").Value = 1.039695968466208
$ws.Range("J6").Value = 1.02702109334905
$ws.Range("K6").Value = 1.026535697595033
$ws.Range("L6").Value = 1.034291444486444
$ws.Range("M6").Value = 1.042110377614219
$ws.Range("N6").Value = 1.013004056116885

# Row 7
$ws.Range("C7").Value = 1.022188748665937
$ws.Range("D7").Value = 1.023777122089088
$ws.Range("E7").Value = 1.031534275845472
$ws.Range("F7").Value = 1.039318179757126
$ws.Range("J7").Value = 1.026785472224582
$ws.Range("K7").Value = 1.02629169510493
$ws.Range("L7").Value = 1.034028892221837
$ws.Range("M7").Value = 1.04179309334369
$ws.Range("N7").Value = 1.012923447932227

# Row 8
$ws.Range("C8").Value = 1.020725911381005
$ws.Range("D8").Value = 1.022502387893733
$ws.Range("E8").Value = 1.030184114099451
$ws.Range("F8").Value = 1.037741509306101
$ws.Range("J8").Value = 1.025800710116415
$ws.Range("K8").Value = 1.02527216389456
$ws.Range("L8").Value = 1.032931983170498
$ws.Range("M8").Value = 1.040468169287144
$ws.Range("N8").Value = 1.012586083693492

# Row 9
$ws.Range("C9").Value = 1.018152436874262
$ws.Range("D9").Value = 1.020262189021229
$ws.Range("E9").Value = 1.027811508835152
$ws.Range("F9").Value = 1.034972331664744
$ws.Range("J9").Value = 1.024065790219161
$ws.Range("K9").Value = 1.023476974988261
$ws.Range("L9").Value = 1.03100103196237
$ws.Range("M9").Value = 1.038138265477069
$ws.Range("N9").Value = 1.011989974284119

# Row 10
$ws.Range("C10").Value = 1.016439591257138
$ws.Range("D10").Value = 1.018772792231163
$ws.Range("E10").Value = 1.026234173601192
$ws.Range("F10").Value = 1.033132319420463
$ws.Range("J10").Value = 1.022909411080321
$ws.Range("K10").Value = 1.022281092265943
$ws.Range("L10").Value = 1.029715031074178
$ws.Range("M10").Value = 1.036588195265998
$ws.Range("N10").Value = 1.011591480050356

# Row 11
$ws.Range("C11").Value = 1.015698565473324
$ws.Range("D11").Value = 1.018128832288675
$ws.Range("E11").Value = 1.025552212411625
$ws.Range("F11").Value = 1.032337016077355
$ws.Range("J11").Value = 1.022408739773229
$ws.Range("K11").Value = 1.021763476684016
$ws.Range("L11").Value = 1.029158484755442
$ws.Range("M11").Value = 1.035917752253727
$ws.Range("N11").Value = 1.011418671126738

# Row 12
$ws.Range("C12").Value = 1.01542341187873
$ws.Range("D12").Value = 1.01788978085279
$ws.Range("E12").Value = 1.02529905698684
$ws.Range("F12").Value = 1.032041819855413
$ws.Range("J12").Value = 1.022222775200686
$ws.Range("K12").Value = 1.021571242525713
$ws.Range("L12").Value = 1.028951803568308
$ws.Range("M12").Value = 1.035668831736868
$ws.Range("N12").Value = 1.0113544434857

# Row 13
$ws.Range("C13").Value = 1.015482428912444
$ws.Range("D13").Value = 1.017941051683968
$ws.Range("E13").Value = 1.025353352659628
$ws.Range("F13").Value = 1.03210513073432
$ws.Range("J13").Value = 1.022262664917883
$ws.Range("K13").Value = 1.021612475990709
$ws.Range("L13").Value = 1.028996135357016
$ws.Range("M13").Value = 1.035722220977581
$ws.Range("N13").Value = 1.011368222283054

# Row 14
$ws.Range("C14").Value = 1.015675819230671
$ws.Range("D14").Value = 1.018109069286104
$ws.Range("E14").Value = 1.02553128333795
$ws.Range("F14").Value = 1.032312610686158
$ws.Range("J14").Value = 1.022393367741351
$ws.Range("K14").Value = 1.021747585910623
$ws.Range("L14").Value = 1.029141399509654
$ws.Range("M14").Value = 1.035897174126639
$ws.Range("N14").Value = 1.011413362839077

# Row 15
$ws.Range("C15").Value = 1.015794986131469
$ws.Range("D15").Value = 1.018212609549888
$ws.Range("E15").Value = 1.025640932876699
$ws.Range("F15").Value = 1.032440474383947
$ws.Range("J15").Value = 1.022473898990778
$ws.Range("K15").Value = 1.02183083573628
$ws.Range("L15").Value = 1.029230907489905
$ws.Range("M15").Value = 1.036004983346332
$ws.Range("N15").Value = 1.011441170298477

# Row 16
$ws.Range("C16").Value = 1.016488784264059
$ws.Range("D16").Value = 1.018815549923921
$ws.Range("E16").Value = 1.026279454913947
$ws.Range("F16").Value = 1.033185131246976
$ws.Range("J16").Value = 1.02294263998695
$ws.Range("K16").Value = 1.022315449107855
$ws.Range("L16").Value = 1.029751973530508
$ws.Range("M16").Value = 1.036632706073421
$ws.Range("N16").Value = 1.011602943375499

# Row 17
$ws.Range("C17").Value = 1.016924158244312
$ws.Range("D17").Value = 1.019194015107871
$ws.Range("E17").Value = 1.026680259955463
$ws.Range("F17").Value = 1.033652618342809
$ws.Range("J17").Value = 1.023236681640314
$ws.Range("K17").Value = 1.022619490252655
$ws.Range("L17").Value = 1.030078904773804
$ws.Range("M17").Value = 1.037026659793874
$ws.Range("N17").Value = 1.011704350266305

# Row 18
$ws.Range("C18").Value = 1.017178167073867
$ws.Range("D18").Value = 1.019414859951583
$ws.Range("E18").Value = 1.026914142668853
$ws.Range("F18").Value = 1.033925434066817
$ws.Range("J18").Value = 1.023408195856578
$ws.Range("K18").Value = 1.022796852535647
$ws.Range("L18").Value = 1.030269627252778
$ws.Range("M18").Value = 1.037256518498098
$ws.Range("N18").Value = 1.011763474262471

# Row 19
$ws.Range("C19").Value = 1.01726478815901
$ws.Range("D19").Value = 1.019490178028921
$ws.Range("E19").Value = 1.026993907584169
$ws.Range("F19").Value = 1.034018480673861
$ws.Range("J19").Value = 1.023466678575723
$ws.Range("K19").Value = 1.022857331926621
$ws.Range("L19").Value = 1.030334663618961
$ws.Range("M19").Value = 1.037334906630304
$ws.Range("N19").Value = 1.011783629786303

# Row 20
$ws.Range("C20").Value = 1.016877440273005
$ws.Range("D20").Value = 1.0191533998076
$ws.Range("E20").Value = 1.02663724703468
$ws.Range("F20").Value = 1.033602447089505
$ws.Range("J20").Value = 1.023205133269757
$ws.Range("K20").Value = 1.022586867424132
$ws.Range("L20").Value = 1.030043825145954
$ws.Range("M20").Value = 1.036984384826868
$ws.Range("N20").Value = 1.011693472844304

# Row 21
$ws.Range("C21").Value = 1.015618867949066
$ws.Range("D21").Value = 1.018059588316606
$ws.Range("E21").Value = 1.025478882902058
$ws.Range("F21").Value = 1.032251507085556
$ws.Range("J21").Value = 1.022354878825014
$ws.Range("K21").Value = 1.021707798543563
$ws.Range("L21").Value = 1.029098621600024
$ws.Range("M21").Value = 1.03584565170029
$ws.Range("N21").Value = 1.011400071138333

# Row 22
$ws.Range("C22").Value = 1.014828110722075
$ws.Range("D22").Value = 1.017372698487876
$ws.Range("E22").Value = 1.024751471951726
$ws.Range("F22").Value = 1.031403360250183
$ws.Range("J22").Value = 1.021820330742941
$ws.Range("K22").Value = 1.021155274152924
$ws.Range("L22").Value = 1.028504594244827
$ws.Range("M22").Value = 1.035130332009973
$ws.Range("N22").Value = 1.011215373908602

# Row 23
$ws.Range("C23").Value = 1.015247253587172
$ws.Range("D23").Value = 1.017736752859941
$ws.Range("E23").Value = 1.025137001161882
$ws.Range("F23").Value = 1.031852861336366
$ws.Range("J23").Value = 1.022103701023302
$ws.Range("K23").Value = 1.02144816067351
$ws.Range("L23").Value = 1.028819474860521
$ws.Range("M23").Value = 1.035509475331173
$ws.Range("N23").Value = 1.011313306599714

# Row 24
$ws.Range("C24").Value = 1.016898549919705
$ws.Range("D24").Value = 1.019171751827567
$ws.Range("E24").Value = 1.026656682413227
$ws.Range("F24").Value = 1.033625116891363
$ws.Range("J24").Value = 1.02321938860577
$ws.Range("K24").Value = 1.02260160821347
$ws.Range("L24").Value = 1.030059676029977
$ws.Range("M24").Value = 1.037003486841567
$ws.Range("N24").Value = 1.011698387960088

# Row 25
$ws.Range("C25").Value = 1.018817244504749
$ws.Range("D25").Value = 1.020840617249542
$ws.Range("E25").Value = 1.028424108096125
$ws.Range("F25").Value = 1.035687152496668
$ws.Range("J25").Value = 1.024514266927051
$ws.Range("K25").Value = 1.02394091389106
$ws.Range("L25").Value = 1.012144275226148
$ws.Range("M25").Value = 1.038740036469494
$ws.Range("N25").Value = 1.012144275225781
